$wb = $excel.ActiveWorkbook

# --- Shared string text update: "Semi Final Points" -> "Workout 6 Points" ---
$wsFM = $wb.Worksheets.Item("FM")
$wsFF = $wb.Worksheets.Item("FF")
$wsFM.Range("C1").Value = "Workout 6 Points"
$wsFF.Range("C1").Value = "Workout 6 Points"

# --- FM (sheet5) new rows 2-4 ---
$wsFM.Range("A2").Value = "Anders Magnus Nes og Anders Vinnes Jacobsen"
$wsFM.Range("B2").Value = 9
$wsFM.Range("C2").Value = 15
$wsFM.Range("D2").Value = 12
$wsFM.Range("E2").Value = 40
$wsFM.Range("F2").Value = 300

$wsFM.Range("A3").Value = "Anders J. Svalestuen og Gabriel Kristiansen"
$wsFM.Range("B3").Value = 9
$wsFM.Range("C3").Value = 10
$wsFM.Range("D3").Value = 13
$wsFM.Range("E3").Value = 4
$wsFM.Range("F3").Value = 300

$wsFM.Range("A4").Value = "Vegard Austrheim Vågen og Henrik Eliassen"
$wsFM.Range("B4").Value = 8
$wsFM.Range("C4").Value = 12
$wsFM.Range("D4").Value = 14
$wsFM.Range("E4").Value = 12
$wsFM.Range("F4").Value = 300

# --- FF (sheet6) new rows 2-4 ---
$wsFF.Range("A2").Value = "Beata Wilman og Ingrid Hamnes"
$wsFF.Range("B2").Value = 8
$wsFF.Range("C2").Value = 15
$wsFF.Range("D2").Value = 12
$wsFF.Range("E2").Value = 45
$wsFF.Range("F2").Value = 300

$wsFF.Range("A3").Value = "Maria Hanssen og Cecilie Rabben"
$wsFF.Range("B3").Value = 9
$wsFF.Range("C3").Value = 12
$wsFF.Range("D3").Value = 13
$wsFF.Range("E3").Value = 8
$wsFF.Range("F3").Value = 300

$wsFF.Range("A4").Value = "Marianne U. Henriksen og Mari S. Andersen"
$wsFF.Range("B4").Value = 9
$wsFF.Range("C4").Value = 12
$wsFF.Range("D4").Value = 13
$wsFF.Range("E4").Value = 21
$wsFF.Range("F4").Value = 300

# --- Column A widths: re-fit now that longer team names were added ---
$wsFM.Columns.Item(1).AutoFit()
$wsFF.Columns.Item(1).AutoFit()

# --- Selections ---
$wsSFM = $wb.Worksheets.Item("SFM")
$wsSFF = $wb.Worksheets.Item("SFF")
[void]$wsSFM.Range("A2:A7").Select()
[void]$wsSFF.Range("A2:A7").Select()
[void]$wsFM.Range("F5").Select()
[void]$wsFF.Range("E9").Select()

# --- Active sheet / tab (FF ends up as the active/selected tab) ---
[void]$wsFF.Activate()

Write-Output "done"
